$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 189, shifting existing rows (189-197) down to (190-198).
$ws.Rows.Item(189).Insert()

# Populate the newly inserted row 189 with the new weekly price record.
$ws.Range("A189").Value = 5
$ws.Range("B189").Value = "Macroferia Regional de Talca"
$ws.Range("C189").Value = "Maule"
$ws.Range("D189").Value = 45147
$ws.Range("E189").Value = 7
$ws.Range("F189").Value = 100112001
$ws.Range("G189").Value = "Berenjena"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 150
$ws.Range("K189").Value = 10000
$ws.Range("L189").Value = 10000
$ws.Range("M189").Value = 10000
$ws.Range("N189").Value = "$/caja 50 unidades"
$ws.Range("O189").Value = "Región de Arica y Parinacota"
$ws.Range("P189").Value = 200
$ws.Range("Q189").Value = 50
$ws.Range("R189").Value = "Hortaliza"
